{"js": "const body = context.document.body;\n\n// The `_GoBack` bookmark (Word's \"last edit location\" marker) is being\n// relocated from the end of the certificate section to the new\n// \"Tool Guides\" heading appended at the end of the document, so the\n// paragraph that used to hold only the bookmark becomes a plain empty\n// paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- Append the new \"Tool Guides\" section at the end of the document ---\nconst toolGuides = body.insertParagraph(\"Tool Guides\", Word.InsertLocation.end);\ntoolGuides.style = \"Heading 1\";\nawait context.sync();\n\n// Re-create the `_GoBack` bookmark right after the new heading's text.\nconst toolGuidesEnd = toolGuides.getRange(Word.RangeLocation.end);\ntoolGuidesEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n\nconst blank1 = body.insertParagraph(\"\", Word.InsertLocation.end);\nblank1.style = \"Normal\";\nawait context.sync();\n\nconst git = body.insertParagraph(\"Git\", Word.InsertLocation.end);\ngit.style = \"Heading 2\";\nawait context.sync();\n\nconst blank2 = body.insertParagraph(\"\", Word.InsertLocation.end);\nblank2.style = \"Normal\";\nawait context.sync();\n\nconst cloning = body.insertParagraph(\"Cloning your own repo on github\", Word.InsertLocation.end);\ncloning.style = \"Heading 3\";\nawait context.sync();\n\nconst codeLines = [\n  \"git clone https://github.com/userName/Repo New_Repo\",\n  \"cd New_Repo\",\n  \"git remote set-url origin https://github.com/userName/New_Repo\",\n  \"git remote add upstream https://github.com/userName/Repo\",\n  \"git push origin master\"\n];\n\nfor (const line of codeLines) {\n  const codePar = body.insertParagraph(line, Word.InsertLocation.end);\n  codePar.style = \"Code\";\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# The `_GoBack` bookmark (Word's \"last edit location\" marker) is being\n# relocated from the end of the certificate section to the new\n# \"Tool Guides\" heading appended at the end of the document, so the\n# paragraph that used to hold only the bookmark becomes a plain empty\n# paragraph.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- Append the new \"Tool Guides\" section at the end of the document ---\n$end = $d.Content\n$end.Collapse(0)\n$end.InsertParagraphAfter()\n$end.Collapse(0)\n$end.Text = \"Tool Guides\"\n$toolGuidesPar = $d.Paragraphs.Last\n$toolGuidesPar.Style = \"Heading 1\"\n\n# Re-create the `_GoBack` bookmark right after the new heading's text.\n$tgRange = $toolGuidesPar.Range\n$bmRange = $d.Range($tgRange.Start, $tgRange.End - 1)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n\n$end = $d.Content\n$end.Collapse(0)\n$end.InsertParagraphAfter()\n$end.Collapse(0)\n$d.Paragraphs.Last.Style = \"Normal\"\n\n$end = $d.Content\n$end.Collapse(0)\n$end.InsertParagraphAfter()\n$end.Collapse(0)\n$end.Text = \"Git\"\n$d.Paragraphs.Last.Style = \"Heading 2\"\n\n$end = $d.Content\n$end.Collapse(0)\n$end.InsertParagraphAfter()\n$end.Collapse(0)\n$d.Paragraphs.Last.Style = \"Normal\"\n\n$end = $d.Content\n$end.Collapse(0)\n$end.InsertParagraphAfter()\n$end.Collapse(0)\n$end.Text = \"Cloning your own repo on github\"\n$d.Paragraphs.Last.Style = \"Heading 3\"\n\n$codeLines = @(\n  \"git clone https://github.com/userName/Repo New_Repo\",\n  \"cd New_Repo\",\n  \"git remote set-url origin https://github.com/userName/New_Repo\",\n  \"git remote add upstream https://github.com/userName/Repo\",\n  \"git push origin master\"\n)\n\nforeach ($line in $codeLines) {\n  $end = $d.Content\n  $end.Collapse(0)\n  $end.InsertParagraphAfter()\n  $end.Collapse(0)\n  $end.Text = $line\n  $d.Paragraphs.Last.Style = \"Code\"\n}\n"}
